$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26..197 down to 27..198.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new data record.
$ws.Range("A26").Value = 8
$ws.Range("B26").Value = "Terminal La Palmera de La Serena"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44473
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 100114013
$ws.Range("G26").Value = "Zanahoria"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 600
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 6500
$ws.Range("N26").Value = "`$/saco 20 kilos"
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 325
$ws.Range("Q26").Value = 20
$ws.Range("R26").Value = "Hortaliza"
